$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 currently has G22/H22 present but empty, and is missing I22/J22.
# Bring it in line with the pattern used by the other data rows (e.g. row 20):
#   - G22, H22 keep their existing style, just need a value of 5
#   - I22, J22 are new cells that need the "total column" style (same as I20:J20) and a value of 5
$ws.Range("I20:J20").Copy() | Out-Null
$ws.Range("I22:J22").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 5
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 5

# Update the active selection in the frozen bottom-right pane to H22.
[void]$ws.Range("H22").Select()
